# "[master] added lab 12" -- the recap/outline slide that used to close out
# this deck (slide 12: "Outline" / "Recap" / "Linking state management with
# networking" / "Exercise") now lives in its own lab-12 deck, so remove it
# from this presentation.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$s.Delete()
